$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# (e.g. "1.001", "0.9999") are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.740.60'

# Row 3
$ws.Range("D3").Value = '1.847.42'
$ws.Range("E3").Value = '  +0.07%  '

# Row 4
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '313.68'
$ws.Range("E5").Value = '  -0.40%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.01%  '

# Row 7
$ws.Range("D7").Value = '0.4311'
$ws.Range("E7").Value = '  +1.57%  '

# Row 8
$ws.Range("D8").Value = '0.3652'
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").Value = '45.01'
$ws.Range("E9").Value = '  -1.77%  '

# Row 10
$ws.Range("D10").Value = '0.07332'
$ws.Range("E10").Value = '  +0.86%  '

# Row 11
$ws.Range("D11").Value = '0.8762'
$ws.Range("E11").Value = '  -2.64%  '

# Row 12
$ws.Range("E12").Value = '  +0.25%  '

# Row 13
$ws.Range("D13").Value = '1.797.37'
$ws.Range("E13").Value = '  +0.70%  '

# Row 14
$ws.Range("D14").Value = '5.339'
$ws.Range("E14").Value = '  -0.75%  '

# Row 15
$ws.Range("D15").Value = '6.520'
$ws.Range("E15").Value = '  -0.65%  '

# Row 16
$ws.Range("D16").Value = '0.06923'
$ws.Range("E16").Value = '  +1.08%  '

# Row 17
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.02%  '

# Row 18
$ws.Range("D18").Value = '79.93'
$ws.Range("E18").Value = '  +2.21%  '

# Row 19
$ws.Range("D19").Value = '0.000009008'
$ws.Range("E19").Value = '  +1.96%  '

# Row 20
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  -0.02%  '

# Row 21
$ws.Range("E21").Value = '  -0.67%  '

# Row 22
$ws.Range("D22").Value = '27.596.63'
$ws.Range("E22").Value = '  -0.14%  '

# Row 23
$ws.Range("D23").Value = '4.978'
$ws.Range("E23").Value = '  +0.10%  '

# Row 24
$ws.Range("E24").Value = '  -2.08%  '

# Row 25
$ws.Range("D25").Value = '2.029.39'
$ws.Range("E25").Value = '  +0.43%  '

# Row 26
$ws.Range("D26").Value = '1.979'
$ws.Range("E26").Value = '  -3.10%  '

# Row 27
$ws.Range("D27").Value = '155.95'
$ws.Range("E27").Value = '  +1.18%  '

# Row 28
$ws.Range("E28").Value = '  +1.98%  '

# Row 29
$ws.Range("D29").Value = '120.04'
$ws.Range("E29").Value = '  +8.56%  '

# Row 30
$ws.Range("D30").Value = '5.254'
$ws.Range("E30").Value = '  -0.47%  '

# Row 31
$ws.Range("D31").Value = '1.859'
$ws.Range("E31").Value = '  +2.68%  '

# Row 32
$ws.Range("D32").Value = '0.08892'
$ws.Range("E32").Value = '  +0.29%  '

# Row 33
$ws.Range("D33").Value = '0.7535'
$ws.Range("E33").Value = '  -2.27%  '

# Row 34
$ws.Range("D34").Value = '4.541'
$ws.Range("E34").Value = '  -0.33%  '

# Row 35
$ws.Range("D35").Value = '2.962'
$ws.Range("E35").Value = '  -0.29%  '

# Row 36
$ws.Range("D36").Value = '1.123'
$ws.Range("E36").Value = '  +3.28%  '

# Row 37
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  +0.00%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.106'
$ws.Range("E38").Value = '  +0.57%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05415'
$ws.Range("E39").Value = '  +0.00%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01931'
$ws.Range("E40").Value = '  +0.04%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.842'
$ws.Range("E41").Value = '  -1.58%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.5083'
$ws.Range("E42").Value = '  +0.31%  '

# Row 43
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.1657'
$ws.Range("E43").Value = '  +0.75%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.674'
$ws.Range("E44").Value = '  -1.91%  '

# Row 45
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '8.340'
$ws.Range("E45").Value = '  +1.41%  '

# Row 46
$ws.Range("D46").Value = '0.06534'
$ws.Range("E46").Value = '  -1.46%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.26'
$ws.Range("E47").Value = '  -1.28%  '

# Row 48
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.4657'
$ws.Range("E48").Value = '  -1.34%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '104.38'
$ws.Range("E49").Value = '  -0.93%  '

# Row 50
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '0.9999'
$ws.Range("E50").Value = '  -0.03%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.623'
$ws.Range("E51").Value = '  -1.01%  '

# Restore the default "Normal" style on column D so no stray text-format
# style is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
